# #5: insurance, claim, debt, investment done
#
# The "保險" (insurance) and "債務" (debt) sheets had their (bold) header
# row accidentally populated with duplicated data values instead of real
# column labels, and were missing the trailing metadata columns
# (property_category / category / date / legislator_name / legislator_id /
# source_file / index, plus a couple of sheet-specific ones) that every
# other sheet ("土地" / "建物" / "汽車") already carries. This fixes the
# header rows and appends the missing metadata columns + their data, using
# the already-present cells on the "土地" sheet as the source of truth for
# both the label text and the shared "normal" / date / legislator / source
# file values (that sheet already stores them as plain text, so copying
# from it avoids Excel's automatic text->date/number coercion).

$wb = $excel.ActiveWorkbook
$land = $wb.Worksheets.Item("土地")

# =======================================================================
# Sheet "保險" (insurance) -- 4th tab
# =======================================================================
$ins = $wb.Worksheets.Item("保險")

# --- Row 1 (bold header) ------------------------------------------------
# B1:D1 currently hold duplicated data values; turn them into the real
# column names and extend the header with the standard metadata columns.
$land.Range("B1").Copy($ins.Range("C1"))   # name
$land.Range("E1").Copy($ins.Range("D1"))   # owner
$land.Range("I1").Copy($ins.Range("E1"))   # property_category
$land.Range("J1").Copy($ins.Range("F1"))   # category
$land.Range("K1").Copy($ins.Range("G1"))   # date
$land.Range("L1").Copy($ins.Range("H1"))   # legislator_name
$land.Range("M1").Copy($ins.Range("I1"))   # legislator_id
$land.Range("N1").Copy($ins.Range("J1"))   # source_file
$land.Range("O1").Copy($ins.Range("K1"))   # index
$ins.Range("B1").Value = "company"

# --- Row 2 ----------------------------------------------------------------
# B2:D2 already hold the correct data (company/name/owner); append the
# metadata columns.
$land.Range("J2").Copy($ins.Range("F2"))   # normal
$land.Range("K2").Copy($ins.Range("G2"))   # 2012-05-01
$land.Range("L2").Copy($ins.Range("H2"))   # 陳歐珀
$land.Range("M2").Copy($ins.Range("I2"))   # 1753
$land.Range("N2").Copy($ins.Range("J2"))   # tmpe4f31
$ins.Range("E2").Value = "insurance"
$ins.Range("K2").Value = 73

# --- Row 3 -----------------------------------------------------------------
$land.Range("J2").Copy($ins.Range("F3"))   # normal
$land.Range("K2").Copy($ins.Range("G3"))   # 2012-05-01
$land.Range("L2").Copy($ins.Range("H3"))   # 陳歐珀
$land.Range("M2").Copy($ins.Range("I3"))   # 1753
$land.Range("N2").Copy($ins.Range("J3"))   # tmpe4f31
$ins.Range("E3").Value = "insurance"
$ins.Range("K3").Value = 74

# =======================================================================
# Sheet "債務" (debt) -- 5th tab
# =======================================================================
$debt = $wb.Worksheets.Item("債務")

# --- Row 1 (bold header) ------------------------------------------------
# B1:G1 currently hold duplicated data values; turn them into the real
# column names and extend the header with the standard metadata columns.
$land.Range("E1").Copy($debt.Range("D1"))   # owner
$land.Range("Q1").Copy($debt.Range("E1"))   # total
$land.Range("F1").Copy($debt.Range("F1"))   # register_date
$land.Range("G1").Copy($debt.Range("G1"))   # register_reason
$land.Range("I1").Copy($debt.Range("H1"))   # property_category
$land.Range("J1").Copy($debt.Range("I1"))   # category
$land.Range("K1").Copy($debt.Range("J1"))   # date
$land.Range("L1").Copy($debt.Range("K1"))   # legislator_name
$land.Range("M1").Copy($debt.Range("L1"))   # legislator_id
$land.Range("N1").Copy($debt.Range("M1"))   # source_file
$land.Range("O1").Copy($debt.Range("N1"))   # index
$debt.Range("B1").Value = "species"
$debt.Range("C1").Value = "debtor"

# --- Row 2 ----------------------------------------------------------------
# B2:G2 already hold the correct data; append the metadata columns.
$land.Range("J2").Copy($debt.Range("I2"))   # normal
$land.Range("K2").Copy($debt.Range("J2"))   # 2012-05-01
$land.Range("L2").Copy($debt.Range("K2"))   # 陳歐珀
$land.Range("M2").Copy($debt.Range("L2"))   # 1753
$land.Range("N2").Copy($debt.Range("M2"))   # tmpe4f31
$debt.Range("H2").Value = "debt"
$debt.Range("N2").Value = 84
